$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.15163
$ws.Range("H2").Value = 9.454890000000001
$ws.Range("I2").Value = 0.0006291248881010851
$ws.Range("J2").Value = 0.0006291248881010851
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.847798333333333
$ws.Range("N2").Value = 5.543395
$ws.Range("O2").Value = 0.05039680725746681
$ws.Range("P2").Value = 0.05039680725746681
$ws.Range("Q2").Value = 5.823576661283334
$ws.Range("R2").Value = 52.41218995155
$ws.Range("S2").Value = 0.00003170588572650576
$ws.Range("T2").Value = 0.00003170588572650576

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.15163
$ws.Range("H3").Value = 9.454890000000001
$ws.Range("I3").Value = 0.0006291248881010851
$ws.Range("J3").Value = 0.0006291248881010851
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.4798556666666667
$ws.Range("N3").Value = 1.439567
$ws.Range("O3").Value = 0.0130875719001099
$ws.Range("P3").Value = 0.0130875719001099
$ws.Range("Q3").Value = 1.512327514736667
$ws.Range("R3").Value = 13.61094763263
$ws.Range("S3").Value = 0.000008233717207171545
$ws.Range("T3").Value = 0.000008233717207171547

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.15163
$ws.Range("H4").Value = 9.454890000000001
$ws.Range("I4").Value = 0.0006291248881010851
$ws.Range("J4").Value = 0.0006291248881010851
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 34.337334
$ws.Range("N4").Value = 103.012002
$ws.Range("O4").Value = 0.9365156208424232
$ws.Range("P4").Value = 0.9365156208424232
$ws.Range("Q4").Value = 108.21857195442
$ws.Range("R4").Value = 973.96714758978
$ws.Range("S4").Value = 0.0005891852851674077
$ws.Range("T4").Value = 0.0005891852851674077

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 4971.754394666666
$ws.Range("H5").Value = 14915.263184
$ws.Range("I5").Value = 0.9924561027819714
$ws.Range("J5").Value = 0.9924561027819713
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.847798333333333
$ws.Range("N5").Value = 5.543395
$ws.Range("O5").Value = 0.05039680725746681
$ws.Range("P5").Value = 0.05039680725746681
$ws.Range("Q5").Value = 9186.799484207742
$ws.Range("R5").Value = 82681.19535786968
$ws.Range("S5").Value = 0.05001661892339967
$ws.Range("T5").Value = 0.05001661892339967

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 4971.754394666666
$ws.Range("H6").Value = 14915.263184
$ws.Range("I6").Value = 0.9924561027819714
$ws.Range("J6").Value = 0.9924561027819713
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.4798556666666667
$ws.Range("N6").Value = 1.439567
$ws.Range("O6").Value = 0.0130875719001099
$ws.Range("P6").Value = 0.0130875719001099
$ws.Range("Q6").Value = 2385.724519555703
$ws.Range("R6").Value = 21471.52067600133
$ws.Range("S6").Value = 0.01298884060286191
$ws.Range("T6").Value = 0.01298884060286191

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 4971.754394666666
$ws.Range("H7").Value = 14915.263184
$ws.Range("I7").Value = 0.9924561027819714
$ws.Range("J7").Value = 0.9924561027819713
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 34.337334
$ws.Range("N7").Value = 103.012002
$ws.Range("O7").Value = 0.9365156208424232
$ws.Range("P7").Value = 0.9365156208424232
$ws.Range("Q7").Value = 170716.7912156371
$ws.Range("R7").Value = 1536451.120940734
$ws.Range("S7").Value = 0.9294506432557097
$ws.Range("T7").Value = 0.9294506432557096

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 34.63986933333333
$ws.Range("H8").Value = 103.919608
$ws.Range("I8").Value = 0.006914772329927541
$ws.Range("J8").Value = 0.006914772329927542
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.847798333333333
$ws.Range("N8").Value = 5.543395
$ws.Range("O8").Value = 0.05039680725746681
$ws.Range("P8").Value = 0.05039680725746681
$ws.Range("Q8").Value = 64.00749282101778
$ws.Range("R8").Value = 576.06743538916
$ws.Range("S8").Value = 0.000348482448340623
$ws.Range("T8").Value = 0.000348482448340623

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 34.63986933333333
$ws.Range("H9").Value = 103.919608
$ws.Range("I9").Value = 0.006914772329927541
$ws.Range("J9").Value = 0.006914772329927542
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.4798556666666667
$ws.Range("N9").Value = 1.439567
$ws.Range("O9").Value = 0.0130875719001099
$ws.Range("P9").Value = 0.0130875719001099
$ws.Range("Q9").Value = 16.62213759219289
$ws.Range("R9").Value = 149.599238329736
$ws.Range("S9").Value = 0.00009049758004081715
$ws.Range("T9").Value = 0.00009049758004081717

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 34.63986933333333
$ws.Range("H10").Value = 103.919608
$ws.Range("I10").Value = 0.006914772329927541
$ws.Range("J10").Value = 0.006914772329927542
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 34.337334
$ws.Range("N10").Value = 103.012002
$ws.Range("O10").Value = 0.9365156208424232
$ws.Range("P10").Value = 0.9365156208424232
$ws.Range("Q10").Value = 1189.440763015024
$ws.Range("R10").Value = 10704.96686713522
$ws.Range("S10").Value = 0.006475792301546101
$ws.Range("T10").Value = 0.006475792301546102

